$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.382.39'
$ws.Range("E2").Value = '  +4.61%  '
$ws.Range("D3").Value = '1.733.98'
$ws.Range("E3").Value = '  +3.28%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '220.62'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +2.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.524'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.20%  '
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '24.42'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  +13.71%  '
$ws.Range("E9").Value = '  +4.57%  '
$ws.Range("E10").Value = '  +2.05%  '
$ws.Range("E11").Value = '  +1.03%  '
$ws.Range("D12").Value = '1.976.32'
$ws.Range("E12").Value = '  +3.19%  '
$ws.Range("D13").Value = '1.735.55'
$ws.Range("E13").Value = '  +2.88%  '
$ws.Range("E14").Value = '  +3.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.563'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  +4.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '67.87'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  +2.40%  '
$ws.Range("D17").Value = '28.323.57'
$ws.Range("E17").Value = '  +4.45%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '243.93'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '8.05'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  -0.31%  '
$ws.Range("D20").Value = '0.0₃0758'
$ws.Range("E20").Value = '  +2.12%  '
$ws.Range("E21").Value = '  -0.16%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.67'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  +3.25%  '
$ws.Range("E23").Value = '  +2.85%  '
$ws.Range("E24").Value = '  +0.65%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '149.41'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  +1.09%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.56'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +4.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.77'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +2.95%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.115'
$ws.Range("D28").ClearFormats()
$ws.Range("E29").Value = '  -0.32%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0514'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +3.25%  '
$ws.Range("E31").Value = '  +3.21%  '
$ws.Range("E32").Value = '  +2.34%  '
$ws.Range("D33").Value = '1.505.81'
$ws.Range("E33").Value = '  -3.96%  '
$ws.Range("E35").Value = '  -1.13%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.971'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +4.04%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.607'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.14%  '
$ws.Range("E38").Value = '  +0.66%  '
$ws.Range("E39").Value = '  +1.44%  '
$ws.Range("E40").Value = '  +1.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '70.89'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.77%  '
$ws.Range("E42").Value = '  +2.74%  '
$ws.Range("E43").Value = '  -0.14%  '
$ws.Range("E44").Value = '  +2.07%  '
$ws.Range("D45").Value = '1.880.57'
$ws.Range("E45").Value = '  +3.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.807'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.05%  '
$ws.Range("E47").Value = '  +9.74%  '
$ws.Range("D48").Value = '0.0₆0115'
$ws.Range("E48").Value = '  +7.56%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '91.22'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +0.57%  '
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '8.28'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +1.71%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.105'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +1.64%  '
